$d = $word.ActiveDocument

function Replace-Text($old, $new) {
  $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                           $true, 1, $false, $new, 2) | Out-Null
}

# 1. Activation date
Replace-Text "Ativação: 01/01/2020" "Ativação: 01/01/2022"

# 2. Objetivos (PT) paragraph - full text replace
$d.Paragraphs(6).Range.Text = "Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental desde a obtenção das amostras in situ até o preparo preliminar do material a ser analisado. Executar procedimentos de análises qualitativas bem como interpretar, avaliar e criticar os resultados obtidos. Objetivos Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Compreender e aplicar os procedimentos mais comuns de amostragem, coleta e preparação de amostras bem como os erros a não cometer no preparo das amostras de interesse ambiental. Empregar tratamentos preliminares no preparo das amostras: limpeza, secagem, moagem e peneiramento. Compreender as bases teóricas da química analítica qualitativa de interesse ambiental."

# 3. Objetivos (EN, italic) paragraph - full text replace
$d.Paragraphs(7).Range.Text = "Train professionals at a higher level with the ability to know the sequence of chemical analysis procedures of environmental interest from obtaining samples in situ to the preliminary preparation of the material to be analyzed. Perform qualitative analysis procedures as well as interpret, evaluate and criticize the results obtained. Specific Objectives: To train professionals at a higher level with the ability to know the steps of the analytical sequence. Understand and apply the most common procedures for sampling, collecting and preparing samples, as well as errors not to make in preparing samples of environmental interest. Use preliminary treatments in sample preparation: cleaning, drying, grinding and sieving. Understand the theoretical bases of qualitative analytical chemistry of environmental interest"

# 4. Docente responsável
Replace-Text "8855158 - Morun Bernardino Neto" "7455355 - Robson da Silva Rocha"

# 5. Programa resumido (PT) paragraph - full text replace
$d.Paragraphs(11).Range.Text = "Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental. Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental."

# 6. Programa resumido (EN, italic) paragraph - full text replace
$d.Paragraphs(12).Range.Text = "Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance."

# 7. Programa (PT) paragraph - full text replace
$d.Paragraphs(14).Range.Text = "- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Análise de sólidos, partículas, sedimentos.- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-).- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes"

# 8. Programa (EN, italic) paragraph - full text replace
$d.Paragraphs(15).Range.Text = "- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples"

# 9. Método text
Replace-Text "O método de avaliação será composto por 2 avaliações teóricas (P1 e P2) de mesmo valor e igual a 10,0 (dez pontos) além de relatórios de atividades de práticas laboratoriais (LT). Cada relatório de atividade laboratorial será avaliado em 10,0 (dez pontos) e a nota final das atividades laboratoriais será obtida pela média aritmética das notas de todos os relatórios." "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais."

# 10. Critério text
Replace-Text "Para o cálculo da nota final será feita a média ponderada das três avaliações descritas acima (P1, P2 e LT) sendo que a avaliação P2 terá peso 2 e as demais terão peso 1, conforme fórmula abaixo:Nota Final=(P_(1 )+ 2P_2+L_T)/4Estará aprovado por notas o aluno que obtiver resultado final igual ou superior a 5,0 pontos." "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."

# 11. Norma de recuperação text
Replace-Text "Entrará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9 (intervalo fechado). Para os alunos em recuperação, a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final, conforme fórmula abaixo:〖Nota Final〗_rec=(Nota Final+P_recuperação)/2" "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado"
